$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1) HashMap sheet: regroup problems.
#    - "LRU Cache" moves here from Divide&Conquer.
#    - New problem "Longest Consecutive Sequence" is added.
#    - The list is re-sorted descending (was ascending).
# ---------------------------------------------------------------------------
$wsHashMap = $wb.Worksheets.Item("HashMap")

# Make room: 8 existing rows -> 10 rows.
$wsHashMap.Rows.Item(9).Insert()
$wsHashMap.Rows.Item(10).Insert()

$wsHashMap.Range("A1").Value = "Valid Sudoku"
$wsHashMap.Range("A2").Value = "Two Sum"
$wsHashMap.Range("A3").Value = "Minimum Window Substring "
$wsHashMap.Range("A4").Value = "Max Points on a Line"
$wsHashMap.Range("A5").Value = "LRU Cache"
$wsHashMap.Range("A6").Value = "Longest Substring Without Repeating Characters"
$wsHashMap.Range("A7").Value = "Longest Consecutive Sequence"
$wsHashMap.Range("A8").Value = "Copy List with Random Pointer"
$wsHashMap.Range("A9").Value = "Anagrams"
$wsHashMap.Range("A10").Value = "4Sum"

# "4Sum" keeps the yellow "done" marker it had (originally on B1, row1),
# now it lives on its new row (B10). B1 keeps its own marker too because it
# never moved (only column A content was reshuffled).
$wsHashMap.Range("B10").Interior.ColorIndex = 6

# Re-sort column A descending (matches the commit's new sort order) so the
# sheet's cached sortState/sortCondition reflect the actual operation used.
$wsHashMap.Sort.SortFields.Clear()
$wsHashMap.Sort.SortFields.Add($wsHashMap.Range("A1"), 0, 2)
$wsHashMap.Sort.SetRange($wsHashMap.Range("A1:A10"))
$wsHashMap.Sort.Header = 0
$wsHashMap.Sort.Apply()

# ---------------------------------------------------------------------------
# 2) Divide&Conquer sheet: "LRU Cache" is removed (it moved to HashMap).
# ---------------------------------------------------------------------------
$wsDivide = $wb.Worksheets.Item("Divide&Conquer")
$wsDivide.Rows.Item(5).Delete()

# ---------------------------------------------------------------------------
# 3) Selections / active sheet & tab bookkeeping.
#    Activate sheets in the order needed so the LAST one activated
#    ("Divide&Conquer") ends up as the workbook's active tab, matching the
#    new activeTab="11" (and no more firstSheet override).
# ---------------------------------------------------------------------------
$wsHashMap.Activate()
$wsHashMap.Range("L20").Select()

$wsDPMatrix = $wb.Worksheets.Item("DP_Matrix")
$wsDPMatrix.Activate()
$wsDPMatrix.Range("A7").Select()

$wsGraphSearch = $wb.Worksheets.Item("Graph&Search")
$wsGraphSearch.Activate()
$wsGraphSearch.Range("C12").Select()

$wsDivide.Activate()
$wsDivide.Range("O16").Select()
